$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    # Force the cell to hold $val as literal text, even when it looks like a
    # number (e.g. "24.47"), matching the source data's inlineStr cells.
    # Briefly switch to a Text number format while assigning, then restore
    # the cell's original style so no visible formatting changes remain.
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

# Rows 2-45: refresh the Price (D) and Volume(1h) (E) columns.
# $null entries mean that column is unchanged for that row.
$deRows = @(
    @(2, "29.332.96", "  -0.02%  "),
    @(3, "1.843.83", "  -0.14%  "),
    @(4, "0.9985", "  +0.09%  "),
    @(5, "240.01", "  -0.05%  "),
    @(6, "0.6269", "  +0.09%  "),
    @(7, "0.9993", "  +0.06%  "),
    @(8, "0.07477", "  -1.50%  "),
    @(9, $null, "  -0.10%  "),
    @(10, "24.47", "  -0.91%  "),
    @(11, "0.07734", "  +0.00%  "),
    @(12, "1.844.04", "  -2.34%  "),
    @(13, "4.981", "  -0.81%  "),
    @(14, "0.6793", "  +0.19%  "),
    @(15, "0.00001051", "  -0.78%  "),
    @(16, "81.95", "  -1.12%  "),
    @(17, $null, "  +0.85%  "),
    @(18, "29.366.63", "  +0.00%  "),
    @(19, "228.79", "  +0.50%  "),
    @(20, "12.31", "  -0.27%  "),
    @(22, "7.494", "  +0.32%  "),
    @(23, "0.9992", "  +0.11%  "),
    @(24, $null, "  +0.02%  "),
    @(25, "8.424", "  +0.02%  "),
    @(26, "0.1366", "  -1.04%  "),
    @(27, "17.50", "  -0.81%  "),
    @(28, "0.06516", "  +16.31%  "),
    @(29, $null, "  -1.83%  "),
    @(30, "1.482", "  +1.68%  "),
    @(31, "4.118", "  +1.32%  "),
    @(32, "4.087", "  -0.36%  "),
    @(33, "1.825", "  -0.17%  "),
    @(34, "1.140", "  -1.73%  "),
    @(35, "0.6936", "  -0.31%  "),
    @(36, $null, "  -0.04%  "),
    @(37, "1.261.29", "  +2.70%  "),
    @(38, "2.838", $null),
    @(39, "0.01834", "  +1.94%  "),
    @(40, "6.780", "  +6.74%  "),
    @(41, "0.9204", "  +2.68%  "),
    @(42, "0.9984", "  +0.00%  "),
    @(43, "2.008.09", "  +1.48%  "),
    @(44, "101.21", "  -0.25%  "),
    @(45, "66.08", "  +0.97%  ")
)

foreach ($row in $deRows) {
    $rownum = $row[0]
    $dval = $row[1]
    $eval = $row[2]
    if ($dval -ne $null) {
        Set-TextValue ($ws.Cells.Item($rownum, 4)) $dval
    }
    if ($eval -ne $null) {
        Set-TextValue ($ws.Cells.Item($rownum, 5)) $eval
    }
}

# Rows 46-51: the coin ranking list shifted - BabyDogeCoin dropped off and
# Cronos newly appeared at the bottom, so every row from 46 to 51 now holds
# different Coin (B), Link (C), Price (D) and Volume(1h) (E) data.
$bcdeRows = @(
    @(46, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "1.726", "  +2.23%  "),
    @(47, "Aptos", "https://coinranking.com/coin/HGYj5JCv5+aptos-apt", "7.071", "  -1.80%  "),
    @(48, "Algorand", "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo", "0.1161", "  +1.83%  "),
    @(49, "EnergySwap", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens", "8.993", "  -0.04%  "),
    @(50, "TheSandbox", "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand", "0.3945", "  -0.99%  "),
    @(51, "Cronos", "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro", "0.05692", "  -0.07%  ")
)

foreach ($row in $bcdeRows) {
    $rownum = $row[0]
    $bval = $row[1]
    $cval = $row[2]
    $dval = $row[3]
    $eval = $row[4]
    Set-TextValue ($ws.Cells.Item($rownum, 2)) $bval
    Set-TextValue ($ws.Cells.Item($rownum, 3)) $cval
    Set-TextValue ($ws.Cells.Item($rownum, 4)) $dval
    Set-TextValue ($ws.Cells.Item($rownum, 5)) $eval
}

"Update complete"
